$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.397.55'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '1.843.53'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9987'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.35'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6341'
$ws.Range('E6').Value = '  +1.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07477'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.17'
$ws.Range('E9').Value = '  +3.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.2906'
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07748'
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('D12').Value = '1.849.80'
$ws.Range('E12').Value = '  +0.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.993'
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6804'
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001024'
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.07'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.274'
$ws.Range('E17').Value = '  +2.80%  '
$ws.Range('D18').Value = '29.374.71'
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '230.77'
$ws.Range('E19').Value = '  +0.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.35'
$ws.Range('E20').Value = '  +0.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9998'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.435'
$ws.Range('E22').Value = '  +0.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9997'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '158.20'
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.513'
$ws.Range('E25').Value = '  +1.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1360'
$ws.Range('E26').Value = '  -1.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.52'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06536'
$ws.Range('E28').Value = '  +14.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.432'
$ws.Range('E29').Value = '  +2.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.487'
$ws.Range('E30').Value = '  +0.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.080'
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.055'
$ws.Range('E32').Value = '  +0.64%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.843'
$ws.Range('E33').Value = '  +1.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.143'
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7004'
$ws.Range('E35').Value = '  +1.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.580'
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('E37').Value = '  +2.58%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '1.251.26'
$ws.Range('E38').Value = '  +0.72%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.818'
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.774'
$ws.Range('E40').Value = '  +4.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9356'
$ws.Range('E41').Value = '  +3.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9999'
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('D43').Value = '2.014.46'
$ws.Range('E43').Value = '  +0.64%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.24'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.49'
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.074'
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.722'
$ws.Range('E47').Value = '  +4.11%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.058'
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1149'
$ws.Range('E49').Value = '  -1.25%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3921'
$ws.Range('E50').Value = '  -0.34%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05677'
$ws.Range('E51').Value = '  -0.09%  '
